$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.370.32'
$ws.Range("E2").Value = '  -0.03%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.567.02'
$ws.Range("E3").Value = '  -0.17%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9978'
$ws.Range("E5").Value = '  -0.34%  '

# Row 6
$ws.Range("E6").Value = '  +0.87%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3765'
$ws.Range("E7").Value = '  +2.56%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.24'
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("E9").Value = '  +0.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07609'
$ws.Range("E10").Value = '  -0.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.141'
$ws.Range("E11").Value = '  -2.37%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.12%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.04'
$ws.Range("E13").Value = '  -0.87%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.977'
$ws.Range("E14").Value = '  -1.34%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.946'
$ws.Range("E15").Value = '  +0.51%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.584.47'
$ws.Range("E16").Value = '  +0.53%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001134'
$ws.Range("E17").Value = '  +0.19%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.92'
$ws.Range("E18").Value = '  +0.36%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06735'
$ws.Range("E19").Value = '  -0.25%  '

# Row 20
$ws.Range("E20").Value = '  +0.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.58'
$ws.Range("E21").Value = '  +0.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.202'
$ws.Range("E22").Value = '  -0.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.94'
$ws.Range("E23").Value = '  -0.30%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.353.72'
$ws.Range("E24").Value = '  -0.11%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.404'
$ws.Range("E25").Value = '  +1.62%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.708'
$ws.Range("E26").Value = '  -6.43%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.13'
$ws.Range("E27").Value = '  +0.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.54'
$ws.Range("E28").Value = '  +1.64%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.023'
$ws.Range("E29").Value = '  +0.95%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.30'
$ws.Range("E30").Value = '  +0.70%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.722.17'
$ws.Range("E31").Value = '  -1.28%  '

# Row 32
$ws.Range("E32").Value = '  -0.09%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.093'
$ws.Range("E33").Value = '  -2.34%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9875'
$ws.Range("E34").Value = '  -4.98%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.14'
$ws.Range("E35").Value = '  -0.53%  '

# Row 36
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08476'
$ws.Range("E36").Value = '  +0.26%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.408'
$ws.Range("E37").Value = '  +7.56%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02519'
$ws.Range("E38").Value = '  -0.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2302'
$ws.Range("E39").Value = '  -0.63%  '

# Row 40
$ws.Range("E40").Value = '  -0.09%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.402'
$ws.Range("E41").Value = '  -2.24%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.36'
$ws.Range("E42").Value = '  -2.89%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6320'
$ws.Range("E43").Value = '  -0.37%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.91'
$ws.Range("E45").Value = '  -1.60%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.811'
$ws.Range("E46").Value = '  +1.59%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5928'
$ws.Range("E47").Value = '  -0.75%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.082'
$ws.Range("E48").Value = '  -0.81%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.273'
$ws.Range("E49").Value = '  +0.80%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.53'
$ws.Range("E50").Value = '  -0.11%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07312'
$ws.Range("E51").Value = '  +0.53%  '
